# Ajuste no formato dos valores de saida das taxas efetivas anuais.
#
# Cria a planilha "Entrada da compra" com o detalhamento da origem dos
# recursos usados na entrada do financiamento, e faz a celula D2 da planilha
# principal somar esses valores via formula em vez de usar um numero fixo.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Cria a nova planilha logo apos "Plan1" ---------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Entrada da compra"

$moneyFormat = '_-"R$"\ * #,##0.00_-;\-"R$"\ * #,##0.00_-;_-"R$"\ * "-"??_-;_-@_-'
$dateFormat  = "mm-dd-yy"

# --- Cabecalho -----------------------------------------------------------------
$ws2.Range("A1").Value = "Fonte"
$ws2.Range("B1").Value = "Data"
$ws2.Range("C1").Value = "Valor Liquido"

# --- Linhas com a origem de cada valor usado na entrada -------------------------
$fontes = @(
    "CDB BANCO MASTER DE INVESTIME - ABR/2030",
    "ARX Hedge FIC INFRA RF",
    "FGTS",
    "Edileuza - Xp CDI Debêntures Incetivadas",
    "Edileuza - ARX Hedge FIC INFRA RF"
)
$valores     = @(109833.27, 116234.39, 219442.32, 11532.87, 47091.43)
$dataEntrada = 45711

for ($i = 0; $i -lt $fontes.Count; $i++) {
    $row = 2 + $i

    $ws2.Range("A$row").Value = $fontes[$i]

    $ws2.Range("B$row").Value2      = $dataEntrada
    $ws2.Range("B$row").NumberFormat = $dateFormat

    $ws2.Range("C$row").Value2      = $valores[$i]
    $ws2.Range("C$row").NumberFormat = $moneyFormat
}

# O nome do fundo "ARX Hedge FIC INFRA RF" tem quebra de texto automatica
$ws2.Range("A3").WrapText = $true

# --- Total (soma dos valores liquidos) ------------------------------------------
$ws2.Range("C7").Formula      = "=SUM(C2:C6)"
$ws2.Range("C7").NumberFormat = $moneyFormat
$ws2.Range("C7").Font.Bold    = $true

# --- Largura das colunas, ajustada ao conteudo ------------------------------------
$ws2.Range("A1").ColumnWidth = 39.94
$ws2.Range("B1").ColumnWidth = 9.72
$ws2.Range("C1").ColumnWidth = 13.17

$ws2.Range("C8").Select() | Out-Null

# --- Planilha principal: entrada passa a vir da soma da nova planilha -----------
$ws1.Range("D2").Formula = "='Entrada da compra'!C7"

$ws1.Activate() | Out-Null
$ws1.Range("D3").Select() | Out-Null
